$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New codes first (so shared-string order matches: Min label, Max label, then descriptions)
$ws.Range("B4").Value = "MiningNumMin"
$ws.Range("B5").Value = "MiningNumMax"

# Row 4: MiningNumMin
$ws.Range("A4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "채굴최소횟수"

# Row 5: MiningNumMax
$ws.Range("A5").Value = 3
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = "채굴최대횟수"

# Apply vertical-center alignment to the new B-column cells (matches style used for MiningNumMin/Max labels)
$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("B4").HorizontalAlignment = 1
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
